# Rework the SignInData sheet: rename/retext the header + data columns,
# swap in the new test e-mail / password values, add a new (empty,
# hyperlink-styled) D2 cell, widen the newly-meaningful columns, and
# refresh the B3:B6 hyperlink's display text to match the new e-mail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value() = "TestCase_Id"
$ws.Range("B1").Value() = "Gmail/LoginEmailId"
$ws.Range("C1").Value() = "Gmail/LoginPassword"

# --- Data rows --------------------------------------------------------
$ws.Range("A2").Value() = "TestCase_01"
$ws.Range("B2").Value() = "ashitestmail@gmail.com"
$ws.Range("C2").Value() = "researcher"

$ws.Range("A3").Value() = "TestCase_02"
$ws.Range("B3").Value() = "ashitestmail@gmail.com"
$ws.Range("C3").Value() = "researcher"

$ws.Range("A4").Value() = "TestCase_03"
$ws.Range("B4").Value() = "ashitestmail@gmail.com"
$ws.Range("C4").Value() = "researcher"

$ws.Range("A5").Value() = "TestCase_04"
$ws.Range("B5").Value() = "ashitestmail@gmail.com"
$ws.Range("C5").Value() = "researcher"

$ws.Range("A6").Value() = "TestCase_05"
$ws.Range("B6").Value() = "ashitestmail@gmail.com"
$ws.Range("C6").Value() = "researcher"

# --- New, still-empty D2 cell, carrying the hyperlink look ------------
$ws.Range("D2").Style = "Hyperlink"

# --- Keep the B3:B6 hyperlink's shown text in sync with the new value -
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $h = $ws.Hyperlinks.Item($i)
    if ($h.Range.Address() -eq '$B$3:$B$6') {
        $h.TextToDisplay = "ashitestmail@gmail.com"
    }
}

# --- Widen the columns that now carry real content ---------------------
$ws.Columns.Item(2).ColumnWidth = 12.7109375
$ws.Columns.Item(3).ColumnWidth = 14.5703125
$ws.Columns.Item(4).ColumnWidth = 15.140625
$ws.Columns.Item(5).ColumnWidth = 17.85546875
